# Update "Horarios" workbook (Linea 141) with the latest scrape results.
# Timestamp of this refresh:
$newTime = "04:49:42"

$wb = $excel.ActiveWorkbook

function Set-Row($ws, [int]$row, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
}

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 16"

# Existing rows 7-12: refresh scrape time + recalculated "Minutos"
$ws1.Range("A7").Value = $newTime
$ws1.Range("D7").Value = 4

$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 27

$ws1.Range("A9").Value = $newTime
$ws1.Range("D9").Value = 33

$ws1.Range("A10").Value = $newTime
$ws1.Range("D10").Value = 45

$ws1.Range("A11").Value = $newTime
$ws1.Range("D11").Value = 57

$ws1.Range("A12").Value = $newTime
$ws1.Range("D12").Value = 65

# New rows 13-21
Set-Row $ws1 13 $newTime "06:04" "16_SANTA ANA" 75 "LP1912"
Set-Row $ws1 14 $newTime "06:11" "215A_EL PATO" 82 "LP1912"
Set-Row $ws1 15 $newTime "06:14" "225_HARAS DEL SUR" 85 "LP1912"
Set-Row $ws1 16 $newTime "06:21" "26_HERNANDEZ" 92 "LP1912"
Set-Row $ws1 17 $newTime "06:27" "23_HERNANDEZ" 98 "LP1912"
Set-Row $ws1 18 $newTime "06:29" "86_EST CHICA-ESC AGRARIA" 100 "LP1912"
Set-Row $ws1 19 $newTime "06:31" "16_SANTA ANA" 102 "LP1912"
Set-Row $ws1 20 $newTime "06:44" "225_C ROCA-H SUR" 115 "LP1912"
Set-Row $ws1 21 $newTime "06:46" "215C_EL PATO" 117 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTime"
$ws2.Range("A3").Value = "Total filas: 4"

# Existing row 7: refresh scrape time + recalculated "Minutos"
$ws2.Range("A7").Value = $newTime
$ws2.Range("D7").Value = 45

# New rows 8-9
Set-Row $ws2 8 $newTime "06:11" "215A_EL PATO" 82 "LP1912"
Set-Row $ws2 9 $newTime "06:46" "215C_EL PATO" 117 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTime"
$ws3.Range("A3").Value = "Total filas: 4"

# A new row is inserted before the former row 6 (which becomes row 7 and is
# otherwise left untouched), then two more rows are appended.
$ws3.Rows.Item(6).Insert()

Set-Row $ws3 6 $newTime "05:43" "215A_LA PLATA" 54 "L6173"
Set-Row $ws3 8 $newTime "06:08" "215A_LA PLATA" 79 "L6173"
Set-Row $ws3 9 $newTime "06:32" "215C_LA PLATA" 103 "L6203"
